# edit.ps1 -- Re-position every shape on slide 1 (the deck's single slide)
# by the same "nudge" the author applied, and apply the two accompanying
# text edits (a new "Hive" line in the Hive/ODBC/JDBC client textbox, and
# splitting "simulated po transactions" into two lines).
#
# PowerPoint's COM Shape.Left/Top are expressed in points (1 pt = 12700 EMU)
# and are backed by a 32-bit float, so naively assigning `emu / 12700.0`
# can land one EMU short after the round-trip through Single precision.
# EmuToPt nudges the point value until it converts back to the exact EMU
# the target OOXML expects.
function EmuToPt {
    param([double]$emu)
    $pt = $emu / 12700.0
    $step = 0.0000001
    $i = 0
    while ($i -lt 20000) {
        $single = [Single]$pt
        $computed = [Math]::Floor([double]$single * 12700.0)
        if ($computed -eq $emu) {
            break
        } elseif ($computed -lt $emu) {
            $pt = $pt + $step
        } else {
            $pt = $pt - $step
        }
        $i = $i + 1
    }
    return $pt
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1 (id=14): Rectangle: Rounded Corners 13
$shape = $s.Shapes.Item(1)
$shape.Left = EmuToPt 3514725
$shape.Top = EmuToPt 3172277

# Shape 2 (id=21): Rectangle: Rounded Corners 20
$shape = $s.Shapes.Item(2)
$shape.Left = EmuToPt 8690696
$shape.Top = EmuToPt 1879888

# Shape 3 (id=22): Rectangle: Rounded Corners 21
$shape = $s.Shapes.Item(3)
$shape.Left = EmuToPt 8703829
$shape.Top = EmuToPt 1192520

# Shape 4 (id=31): Rectangle: Rounded Corners 30
$shape = $s.Shapes.Item(4)
$shape.Left = EmuToPt 8696255
$shape.Top = EmuToPt 394631

# Shape 5 (id=32): Straight Arrow Connector 31
$shape = $s.Shapes.Item(5)
$shape.Left = EmuToPt 9581284
$shape.Top = EmuToPt 970430

# Shape 6 (id=33): Straight Arrow Connector 32
$shape = $s.Shapes.Item(6)
$shape.Left = EmuToPt 9581284
$shape.Top = EmuToPt 1725278

# Shape 7 (id=38): Smiley Face 37
$shape = $s.Shapes.Item(7)
$shape.Left = EmuToPt 10945660
$shape.Top = EmuToPt 470831

# Shape 8 (id=39): Straight Arrow Connector 38
$shape = $s.Shapes.Item(8)
$shape.Left = EmuToPt 10589715
$shape.Top = EmuToPt 627343

# Shape 9 (id=50): TextBox 49
$shape = $s.Shapes.Item(9)
$shape.Left = EmuToPt 10557463
$shape.Top = EmuToPt 834541

# Shape 10 (id=60): Rectangle: Rounded Corners 15
$shape = $s.Shapes.Item(10)
$shape.Left = EmuToPt 217168
$shape.Top = EmuToPt 2853582

# Shape 11 (id=61): Rectangle: Rounded Corners 14
$shape = $s.Shapes.Item(11)
$shape.Left = EmuToPt 290873
$shape.Top = EmuToPt 3259397

# Shape 12 (id=62): Straight Arrow Connector 61
$shape = $s.Shapes.Item(12)
$shape.Left = EmuToPt 2105138
$shape.Top = EmuToPt 3617225

# Shape 13 (id=64): Rectangle: Rounded Corners 14
$shape = $s.Shapes.Item(13)
$shape.Left = EmuToPt 405173
$shape.Top = EmuToPt 3411797

# Shape 14 (id=65): Rectangle: Rounded Corners 14
$shape = $s.Shapes.Item(14)
$shape.Left = EmuToPt 498840
$shape.Top = EmuToPt 3564197

# Shape 15 (id=81): TextBox 80
$shape = $s.Shapes.Item(15)
$shape.Left = EmuToPt 10911463
$shape.Top = EmuToPt 2641022
$tr = $shape.TextFrame.TextRange
$tr.InsertBefore("Hive`r") | Out-Null

# Shape 16 (id=69): Rectangle: Rounded Corners 15
$shape = $s.Shapes.Item(16)
$shape.Left = EmuToPt 3322978
$shape.Top = EmuToPt 4244883

# Shape 17 (id=78): Rectangle: Rounded Corners 14
$shape = $s.Shapes.Item(17)
$shape.Left = EmuToPt 3514201
$shape.Top = EmuToPt 4650698

# Shape 18 (id=82): Rectangle: Rounded Corners 14
$shape = $s.Shapes.Item(18)
$shape.Left = EmuToPt 3666601
$shape.Top = EmuToPt 4803098

# Shape 19 (id=83): Rectangle: Rounded Corners 14
$shape = $s.Shapes.Item(19)
$shape.Left = EmuToPt 3819001
$shape.Top = EmuToPt 4955498

# Shape 20 (id=84): Straight Arrow Connector 83
$shape = $s.Shapes.Item(20)
$shape.Left = EmuToPt 4827446
$shape.Top = EmuToPt 3825351

# Shape 21 (id=85): Straight Arrow Connector 84
$shape = $s.Shapes.Item(21)
$shape.Left = EmuToPt 3796408
$shape.Top = EmuToPt 3793193

# Shape 22 (id=86): Rectangle: Rounded Corners 14
$shape = $s.Shapes.Item(22)
$shape.Left = EmuToPt 1386818
$shape.Top = EmuToPt 4886668

# Shape 23 (id=87): Straight Arrow Connector 86
$shape = $s.Shapes.Item(23)
$shape.Left = EmuToPt 2676048
$shape.Top = EmuToPt 5241679

# Shape 24 (id=2): TextBox 1
$shape = $s.Shapes.Item(24)
$shape.Left = EmuToPt 1777257
$shape.Top = EmuToPt 3283384
$tr = $shape.TextFrame.TextRange
$tr.InsertBefore("simulated `r") | Out-Null
$para2 = $tr.Paragraphs(2)
$lead = $tr.Characters($para2.Start, 10)
$lead.Text = ""

# Shape 25 (id=96): TextBox 95
$shape = $s.Shapes.Item(25)
$shape.Left = EmuToPt 4015116
$shape.Top = EmuToPt 3676859

# Shape 26 (id=97): TextBox 96
$shape = $s.Shapes.Item(26)
$shape.Left = EmuToPt 2948848
$shape.Top = EmuToPt 3717123

# Shape 27 (id=88): Rectangle: Rounded Corners 35
$shape = $s.Shapes.Item(27)
$shape.Left = EmuToPt 5961033
$shape.Top = EmuToPt 1437460

# Shape 28 (id=90): Rectangle: Rounded Corners 27
$shape = $s.Shapes.Item(28)
$shape.Left = EmuToPt 6156554
$shape.Top = EmuToPt 1961518

# Shape 29 (id=99): Straight Arrow Connector 98
$shape = $s.Shapes.Item(29)
$shape.Left = EmuToPt 5191396
$shape.Top = EmuToPt 2262166

# Shape 30 (id=103): Straight Arrow Connector 102
$shape = $s.Shapes.Item(30)
$shape.Left = EmuToPt 8289394
$shape.Top = EmuToPt 2347649

# Shape 31 (id=104): Straight Arrow Connector 103
$shape = $s.Shapes.Item(31)
$shape.Left = EmuToPt 10589715
$shape.Top = EmuToPt 3091236

# Shape 32 (id=100): Rectangle: Rounded Corners 27
$shape = $s.Shapes.Item(32)
$shape.Left = EmuToPt 6308954
$shape.Top = EmuToPt 2113918

# Shape 33 (id=107): Rectangle: Rounded Corners 35
$shape = $s.Shapes.Item(33)
$shape.Left = EmuToPt 6045271
$shape.Top = EmuToPt 4300403

# Shape 34 (id=106): Rectangle: Rounded Corners 27
$shape = $s.Shapes.Item(34)
$shape.Left = EmuToPt 6244743
$shape.Top = EmuToPt 4723739

# Shape 35 (id=108): Rectangle: Rounded Corners 27
$shape = $s.Shapes.Item(35)
$shape.Left = EmuToPt 6397143
$shape.Top = EmuToPt 4876139

# Shape 36 (id=109): Rectangle: Rounded Corners 21
$shape = $s.Shapes.Item(36)
$shape.Left = EmuToPt 8695488
$shape.Top = EmuToPt 4681178

# Shape 37 (id=110): Straight Arrow Connector 109
$shape = $s.Shapes.Item(37)
$shape.Left = EmuToPt 8421883
$shape.Top = EmuToPt 5000876

# Shape 38 (id=111): Rectangle: Rounded Corners 35
$shape = $s.Shapes.Item(38)
$shape.Left = EmuToPt 6012896
$shape.Top = EmuToPt 2862989

# Shape 39 (id=112): Rectangle: Rounded Corners 27
$shape = $s.Shapes.Item(39)
$shape.Left = EmuToPt 6208417
$shape.Top = EmuToPt 3387047

# Shape 40 (id=113): Rectangle: Rounded Corners 27
$shape = $s.Shapes.Item(40)
$shape.Left = EmuToPt 6360817
$shape.Top = EmuToPt 3539447

# Shape 41 (id=114): Straight Arrow Connector 113
$shape = $s.Shapes.Item(41)
$shape.Left = EmuToPt 5354287
$shape.Top = EmuToPt 3438912

# Shape 42 (id=115): Rectangle: Rounded Corners 21
$shape = $s.Shapes.Item(42)
$shape.Left = EmuToPt 8687225
$shape.Top = EmuToPt 3493039

# Shape 43 (id=116): Straight Arrow Connector 115
$shape = $s.Shapes.Item(43)
$shape.Left = EmuToPt 8401878
$shape.Top = EmuToPt 3708809

# Shape 44 (id=117): Rectangle: Rounded Corners 21
$shape = $s.Shapes.Item(44)
$shape.Left = EmuToPt 8684176
$shape.Top = EmuToPt 2831501

# Shape 45 (id=118): Straight Arrow Connector 117
$shape = $s.Shapes.Item(45)
$shape.Left = EmuToPt 9574764
$shape.Top = EmuToPt 2665564

# Shape 46 (id=58): Straight Arrow Connector 57
$shape = $s.Shapes.Item(46)
$shape.Left = EmuToPt 5288907
$shape.Top = EmuToPt 3717077

# Shape 47 (id=120): TextBox 119
$shape = $s.Shapes.Item(47)
$shape.Left = EmuToPt 5730941
$shape.Top = EmuToPt 5842337

# Shape 48 (id=124): Straight Arrow Connector 123
$shape = $s.Shapes.Item(48)
$shape.Left = EmuToPt 6613403
$shape.Top = EmuToPt 5596690
